# Insert a new weekly price record as row 133 on the single sheet,
# shifting the existing rows 133:186 down to 134:187.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(133).Insert()

$ws.Range("A133").Value = 4
$ws.Range("B133").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C133").Value = "Los Lagos"
$ws.Range("D133").Value = 45027
$ws.Range("E133").Value = 10
$ws.Range("F133").Value = 100112052
$ws.Range("G133").Value = "Albahaca"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 90
$ws.Range("K133").Value = 8000
$ws.Range("L133").Value = 8000
$ws.Range("M133").Value = 8000
$ws.Range("N133").Value = "`$/docena de matas"
$ws.Range("O133").Value = "Región Metropolitana"
$ws.Range("P133").Value = 1333
$ws.Range("Q133").Value = 6
$ws.Range("R133").Value = "Hortaliza"
